$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.494.02"
$ws.Range("E2").Value = "  +5.19%  "
$ws.Range("D3").Value = "1.725.27"
$ws.Range("E3").Value = "  +4.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5351"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2660"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06600"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07656"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.606"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "1.734.49"
$ws.Range("E13").Value = "  +5.47%  "
$ws.Range("D14").Value = "1.960.90"
$ws.Range("E14").Value = "  +4.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5795"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.81%  "
$ws.Range("D16").Value = "0.0₅8287"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.99%  "
$ws.Range("D18").Value = "27.469.31"
$ws.Range("E18").Value = "  +5.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +12.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.723"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.021"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.006"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.95%  "
$ws.Range("E26").Value = "  +13.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1231"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.327"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.95%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05468"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.300"
$ws.Range("D31").Style = "Normal"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.436"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.655"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.863"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9584"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.430"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5923"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01643"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.908"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("D41").Value = "1.047.98"
$ws.Range("E41").Value = "  +1.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8463"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.005"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").Value = "1.867.26"
$ws.Range("E45").Value = "  +4.46%  "
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4504"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.163"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05251"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.46%  "
